$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.763.04"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "1.724.69"
$ws.Range("E3").Value = "  -0.04%  "
Set-TextValue $ws.Range("D4") "0.9974"
$ws.Range("E4").Value = "  -0.23%  "
Set-TextValue $ws.Range("D5") "241.06"
$ws.Range("E5").Value = "  -0.93%  "
Set-TextValue $ws.Range("D6") "0.9980"
$ws.Range("E6").Value = "  -0.18%  "
Set-TextValue $ws.Range("D7") "0.4851"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("E8").Value = "  -1.31%  "
Set-TextValue $ws.Range("D9") "0.06198"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "1.727.67"
$ws.Range("E10").Value = "  +0.06%  "
Set-TextValue $ws.Range("D11") "15.97"
$ws.Range("E11").Value = "  +3.30%  "
Set-TextValue $ws.Range("D12") "0.06911"
$ws.Range("E12").Value = "  -1.25%  "
Set-TextValue $ws.Range("D13") "0.6079"
$ws.Range("E13").Value = "  +1.24%  "
Set-TextValue $ws.Range("D14") "4.477"
$ws.Range("E14").Value = "  -1.53%  "
Set-TextValue $ws.Range("D15") "76.86"
$ws.Range("E15").Value = "  -0.85%  "
Set-TextValue $ws.Range("D16") "0.9984"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "26.564.05"
$ws.Range("E17").Value = "  +0.51%  "
Set-TextValue $ws.Range("D18") "0.9973"
$ws.Range("E18").Value = "  -0.23%  "
Set-TextValue $ws.Range("D19") "0.000007157"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "1.950.19"
$ws.Range("E21").Value = "  +0.15%  "
Set-TextValue $ws.Range("D22") "4.429"
$ws.Range("E22").Value = "  -1.13%  "
Set-TextValue $ws.Range("D23") "8.559"
$ws.Range("E23").Value = "  -0.37%  "
Set-TextValue $ws.Range("D24") "5.064"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -0.27%  "
Set-TextValue $ws.Range("D27") "1.771"
$ws.Range("E27").Value = "  +2.93%  "
Set-TextValue $ws.Range("D28") "1.379"
$ws.Range("E28").Value = "  -1.38%  "
Set-TextValue $ws.Range("D29") "105.89"
$ws.Range("E29").Value = "  -0.96%  "
Set-TextValue $ws.Range("D30") "3.935"
$ws.Range("E30").Value = "  -0.45%  "
Set-TextValue $ws.Range("D31") "0.07948"
$ws.Range("E31").Value = "  -0.82%  "
Set-TextValue $ws.Range("D32") "3.694"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D34") "0.9971"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D35") "2.597"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D36") "1.009"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.6211"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "0.9275"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "2.039"
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.433"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D41") "0.9972"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.01495"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "5.643"
$ws.Range("E43").Value = "  +5.51%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D44") "99.66"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D45") "0.3835"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "6.843"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D47") "0.1157"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.05391"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "7.875"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D50") "30.12"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "51.49"
$ws.Range("E51").Value = "  +1.08%  "
